$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the IP address column (D) from the old 33.x.x.x network to the new 172.x.x.x network
$ws.Range("D3").Value = "172.0.0.0"
$ws.Range("D4").Value = "172.18.0.0"
$ws.Range("D5").Value = "172.13.0.0"
$ws.Range("D6").Value = "172.18.29.0"
$ws.Range("D7").Value = "172.18.30.0"
$ws.Range("D8").Value = "172.13.28.0"
$ws.Range("D9").Value = "172.13.20.0"

# Update the selected cell to reflect where the user's cursor ended up
$ws.Range("D9").Select()
